# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets, as produced by the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 121
    $ws.Range("F3").Value = 5187
    $ws.Range("F7").Value = 786

    if ($name -eq "展览") {
        $ws.Range("F8").Value = 276
    } else {
        $ws.Range("F9").Value = 276
    }
}
